$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per the diff: row number -> column letter -> new numeric value.
$data = @{
    2 = @{ C = 1.574017763137817; D = 0; E = 3097.927280284184; F = 0.1604769829084368; G = 0.1588318131711723; H = 0.1453436509026099; I = 0.1192395740986096; J = 0.1043013785545607; K = 0.09380275646426332; L = 0.08476786689685835; M = 0.07898404109961248; N = 0.07517655251836362; O = 0.07209565543819445; P = 0.06909174256330598; Q = 0.0669768684172341; R = 0.06544886804503286; S = 0.06341059408460627; T = 0.06229062839781802; U = 0.06197782615479622; V = 0.06135780827652566; W = 0.06086045663420624; X = 0.06049406518640967; Y = 0.06038844600943827 }
    3 = @{ C = 1.638013124465942; E = 3227.34484578137; F = 0.1579152468888128; G = 0.1412008429573984; H = 0.1412008429573984; I = 0.1384717320572099; J = 0.112865719140762; K = 0.1022186529250497; L = 0.09116065041110696; M = 0.08613685500236777; N = 0.07950514317872426; O = 0.07574116805960003; P = 0.07284303571626514; Q = 0.07062896024317021; R = 0.06836927076045866; S = 0.06651100357673552; T = 0.06564280613724381; U = 0.06483665944385887; V = 0.06426812444042525; W = 0.06315104071201492; X = 0.06315104071201492; Y = 0.06291120557078693 }
    4 = @{ C = 1.642001152038574; E = 3643.767708799376; F = 0.1589803388499471; G = 0.1421050351351029; H = 0.1371123582911028; I = 0.1371123582911028; J = 0.1371123582911028; K = 0.1302823954030578; L = 0.117543306052887; M = 0.1049357001073789; N = 0.09882248057202381; O = 0.09165549893502091; P = 0.08652756311196635; Q = 0.08189060157074844; R = 0.07973481403823296; S = 0.07695229253058808; T = 0.07579953287714897; U = 0.07357336647325916; V = 0.0728850275861682; W = 0.07162538570104197; X = 0.07133149508198004; Y = 0.07102861030798002 }
    5 = @{ C = 1.441996335983276; E = 3041.61555505221; F = 0.160543468901924; G = 0.150623564293921; H = 0.1394917005067958; I = 0.1205068396276678; J = 0.1062974719832357; K = 0.09470706416902007; L = 0.08320693184407207; M = 0.07740338584114539; N = 0.07154422109194554; O = 0.06970215344897202; P = 0.06690221511636493; Q = 0.06591807215124143; R = 0.06354358764565454; S = 0.06238641413368889; T = 0.06146766628371915; U = 0.06066715931053943; V = 0.06029829377244593; W = 0.05990174609430413; X = 0.05940899115773397; Y = 0.05929075156047191 }
    6 = @{ C = 1.935002088546753; E = 2867.279455970356; F = 0.1592364624934904; G = 0.1592364624934904; H = 0.1273642461882175; I = 0.1058743284319098; J = 0.09657306432779135; K = 0.08584860248738614; L = 0.08025652437970295; M = 0.07456015184770641; N = 0.06895595473332161; O = 0.06663919736600236; P = 0.06454416336785292; Q = 0.06167380849407362; R = 0.0604946964057464; S = 0.05886137396995898; T = 0.05810157901524822; U = 0.0576115096282531; V = 0.05685386045265239; W = 0.0565247909798257; X = 0.05625426565374983; Y = 0.05589238705595235 }
    7 = @{ C = 1.598997354507446; E = 2898.013662203069; F = 0.1563813645940376; G = 0.1524910398811586; H = 0.1182023827869908; I = 0.1022060066474225; J = 0.09234209774111667; K = 0.08107478813149183; L = 0.07686804236450431; M = 0.07123555062829283; N = 0.06859133972825146; O = 0.06602517531629666; P = 0.06452052814850333; Q = 0.06224500430263762; R = 0.06101057941329222; S = 0.05968050183373975; T = 0.05876978212625195; U = 0.05806177314272683; V = 0.0576018677114106; W = 0.0569828192974064; X = 0.05660037763520541; Y = 0.05649149438992335 }
    8 = @{ C = 2.003012657165527; E = 2939.709604098449; F = 0.1600589687846479; G = 0.1409645909646496; H = 0.1298292210015428; I = 0.1097505933768559; J = 0.09774246348089675; K = 0.08771252754728166; L = 0.08069766052270609; M = 0.07323832467553289; N = 0.07079435054924682; O = 0.06693874997809555; P = 0.06508327828417342; Q = 0.06350155109263192; R = 0.0618784935929981; S = 0.06031190251786703; T = 0.05975325062583559; U = 0.05892652629556996; V = 0.05821447632848823; W = 0.05779348245606464; X = 0.0574494382531978; Y = 0.05730428078164616 }
    9 = @{ C = 1.504004240036011; E = 3102.212111944378; F = 0.1577898427418709; G = 0.1425492055572619; H = 0.1397772226173139; I = 0.1248980761070074; J = 0.1046411697356955; K = 0.09943109442559167; L = 0.08996647514315891; M = 0.08115496921861161; N = 0.07547582119351029; O = 0.0707706391583217; P = 0.06855657012814691; Q = 0.06656993727560336; R = 0.06545312250463932; S = 0.06383188326697257; T = 0.06297585940775081; U = 0.06204000552490244; V = 0.06149754954465499; W = 0.06087802853453047; X = 0.06063788709515365; Y = 0.06047197099306778 }
    10 = @{ C = 1.452996969223022; E = 3012.736631736908; F = 0.1561480874227308; G = 0.1439325399783033; H = 0.1316571244770818; I = 0.1148656027738071; J = 0.1031897377690847; K = 0.09186793517604788; L = 0.08578611892288779; M = 0.07973066259662318; N = 0.07308808984826293; O = 0.06929143909120569; P = 0.06573715888735784; Q = 0.06422260995744615; R = 0.06319127768606801; S = 0.06199246366727235; T = 0.06133623535375234; U = 0.06071831863642067; V = 0.05983181511877048; W = 0.05939433981745723; X = 0.05898353585599979; Y = 0.05872780958551475 }
    11 = @{ C = 1.624996900558472; E = 2930.416218952965; F = 0.1592694089928404; G = 0.1581492342732628; H = 0.1331072492332419; I = 0.1122114271957293; J = 0.09073330689490262; K = 0.08838555353650333; L = 0.08021136356524923; M = 0.07380946005141852; N = 0.0704479678855874; O = 0.06717342662143357; P = 0.06559506911034349; Q = 0.06249057542791098; R = 0.06161282296562198; S = 0.05985961337896164; T = 0.05918948791327575; U = 0.0588870254255883; V = 0.05801009950308132; W = 0.05761128231723563; X = 0.0573432678508763; Y = 0.05712312317647104 }
}

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($col in $rowData.Keys) {
        $colIdx = [int][char]$col - [int][char]'A' + 1
        $ws.Cells.Item([int]$r, $colIdx).Value = $rowData[$col]
    }
}
